$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 26
$lastColBefore = 32   # A..AF (32 columns) before the edit

# Remember the current (about-to-be-overwritten) column B contents -- the
# latest price snapshot for every SKU row -- before the shift begins.
$colBValues = @{}
for ($row = 1; $row -le $lastRow; $row++) {
    $colBValues[$row] = $ws.Cells.Item($row, 2).Value2
}

# Shift every price-history column (B..AF) one place to the right (C..AG),
# working from the rightmost column down to B so values are not clobbered
# before they are read.
for ($col = $lastColBefore; $col -ge 2; $col--) {
    for ($row = 1; $row -le $lastRow; $row++) {
        $srcVal = $ws.Cells.Item($row, $col).Value2
        $ws.Cells.Item($row, $col + 1).Value = $srcVal
    }
}

# New column B becomes the most recent snapshot: a fresh timestamp header,
# with each row's price carried over unchanged from the prior latest column.
$ws.Cells.Item(1, 2).Value = "2025-12-23 10:31"
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 2).Value = $colBValues[$row]
}

# Give the freshly used last column (AG, index 33) the same display width
# as every other price column so the sheet stays visually consistent.
$ws.Columns.Item($lastColBefore + 1).ColumnWidth = 20.17
